# Update cryptos list (prices & 1h volume %) to reflect the latest scrape,
# and shift rows 46-51 up by one (BabyDogeCoin row dropped, Cronos appended).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-45: Price (D) and Volume(1h) (E) updates; Coin/Link unchanged ---
$ws.Cells.Item(2, 4).Value = "29.341.15"
$ws.Cells.Item(2, 5).Value = "  -0.11%  "
$ws.Cells.Item(3, 4).Value = "1.840.13"
$ws.Cells.Item(3, 5).Value = "  -0.29%  "
$ws.Cells.Item(4, 4).Value = "'0.9990"
$ws.Cells.Item(4, 5).Value = "  +0.10%  "
$ws.Cells.Item(5, 4).Value = "'239.26"
$ws.Cells.Item(5, 5).Value = "  -0.54%  "
$ws.Cells.Item(6, 5).Value = "  -0.32%  "
$ws.Cells.Item(7, 5).Value = "  +0.10%  "
$ws.Cells.Item(8, 4).Value = "'0.07432"
$ws.Cells.Item(8, 5).Value = "  -0.94%  "
$ws.Cells.Item(9, 5).Value = "  +2.23%  "
$ws.Cells.Item(10, 4).Value = "'0.2891"
$ws.Cells.Item(10, 5).Value = "  -0.59%  "
$ws.Cells.Item(11, 4).Value = "'0.07731"
$ws.Cells.Item(11, 5).Value = "  +0.24%  "
$ws.Cells.Item(12, 4).Value = "1.838.85"
$ws.Cells.Item(12, 5).Value = "  -0.34%  "
$ws.Cells.Item(13, 4).Value = "'4.957"
$ws.Cells.Item(13, 5).Value = "  -1.00%  "
$ws.Cells.Item(14, 4).Value = "'0.6740"
$ws.Cells.Item(14, 5).Value = "  -0.99%  "
$ws.Cells.Item(15, 4).Value = "'0.00001020"
$ws.Cells.Item(15, 5).Value = "  -1.09%  "
$ws.Cells.Item(16, 4).Value = "'81.56"
$ws.Cells.Item(16, 5).Value = "  -0.79%  "
$ws.Cells.Item(17, 4).Value = "'6.218"
$ws.Cells.Item(17, 5).Value = "  +1.01%  "
$ws.Cells.Item(18, 4).Value = "29.348.65"
$ws.Cells.Item(18, 5).Value = "  -0.11%  "
$ws.Cells.Item(19, 4).Value = "'228.98"
$ws.Cells.Item(19, 5).Value = "  -0.18%  "
$ws.Cells.Item(20, 5).Value = "  -0.52%  "
$ws.Cells.Item(21, 5).Value = "  +0.11%  "
$ws.Cells.Item(22, 4).Value = "'7.331"
$ws.Cells.Item(22, 5).Value = "  -1.61%  "
$ws.Cells.Item(23, 4).Value = "'1.001"
$ws.Cells.Item(23, 5).Value = "  +0.21%  "
$ws.Cells.Item(24, 4).Value = "'158.05"
$ws.Cells.Item(24, 5).Value = "  -0.52%  "
$ws.Cells.Item(25, 4).Value = "'8.464"
$ws.Cells.Item(25, 5).Value = "  +0.59%  "
$ws.Cells.Item(26, 4).Value = "'0.1345"
$ws.Cells.Item(26, 5).Value = "  -2.55%  "
$ws.Cells.Item(27, 4).Value = "'17.36"
$ws.Cells.Item(27, 5).Value = "  -1.13%  "
$ws.Cells.Item(28, 4).Value = "'0.07407"
$ws.Cells.Item(28, 5).Value = "  +15.48%  "
$ws.Cells.Item(29, 4).Value = "'1.463"
$ws.Cells.Item(29, 5).Value = "  +5.51%  "
$ws.Cells.Item(30, 4).Value = "'1.477"
$ws.Cells.Item(30, 5).Value = "  +0.28%  "
$ws.Cells.Item(31, 4).Value = "'4.035"
$ws.Cells.Item(31, 5).Value = "  -1.42%  "
$ws.Cells.Item(32, 4).Value = "'4.035"
$ws.Cells.Item(32, 5).Value = "  -0.48%  "
$ws.Cells.Item(33, 4).Value = "'1.819"
$ws.Cells.Item(33, 5).Value = "  -0.02%  "
$ws.Cells.Item(34, 4).Value = "'1.138"
$ws.Cells.Item(34, 5).Value = "  -0.37%  "
$ws.Cells.Item(35, 4).Value = "'0.6943"
$ws.Cells.Item(35, 5).Value = "  -0.65%  "
$ws.Cells.Item(37, 4).Value = "'0.01837"
$ws.Cells.Item(37, 5).Value = "  +0.85%  "
$ws.Cells.Item(38, 4).Value = "'2.801"
$ws.Cells.Item(38, 5).Value = "  -1.02%  "
$ws.Cells.Item(39, 4).Value = "'6.844"
$ws.Cells.Item(39, 5).Value = "  +3.94%  "
$ws.Cells.Item(40, 4).Value = "1.231.48"
$ws.Cells.Item(40, 5).Value = "  -2.18%  "
$ws.Cells.Item(41, 4).Value = "'0.9318"
$ws.Cells.Item(41, 5).Value = "  +2.67%  "
$ws.Cells.Item(42, 4).Value = "'0.9999"
$ws.Cells.Item(42, 5).Value = "  +0.15%  "
$ws.Cells.Item(43, 4).Value = "1.980.83"
$ws.Cells.Item(43, 5).Value = "  -1.26%  "
$ws.Cells.Item(44, 4).Value = "'100.53"
$ws.Cells.Item(44, 5).Value = "  -0.80%  "
$ws.Cells.Item(45, 4).Value = "'65.25"
$ws.Cells.Item(45, 5).Value = "  -1.61%  "

# --- Rows 46-51: BabyDogeCoin row removed -> rows shift up, Cronos appended at the end ---
$ws.Cells.Item(46, 2).Value = "RenderToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(46, 4).Value = "'1.703"
$ws.Cells.Item(46, 5).Value = "  -0.20%  "
$ws.Cells.Item(47, 2).Value = "Aptos"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(47, 4).Value = "'6.939"
$ws.Cells.Item(47, 5).Value = "  -1.64%  "
$ws.Cells.Item(48, 2).Value = "EnergySwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(48, 4).Value = "'8.905"
$ws.Cells.Item(48, 5).Value = "  -1.31%  "
$ws.Cells.Item(49, 2).Value = "Algorand"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(49, 4).Value = "'0.1138"
$ws.Cells.Item(49, 5).Value = "  -3.28%  "
$ws.Cells.Item(50, 2).Value = "TheSandbox"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(50, 4).Value = "'0.3901"
$ws.Cells.Item(50, 5).Value = "  -1.04%  "
$ws.Cells.Item(51, 2).Value = "Cronos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(51, 4).Value = "'0.05668"
$ws.Cells.Item(51, 5).Value = "  -0.75%  "
